# Insert a new first column ("owner_name") into the StockInstrument sheet.
# This mirrors a user selecting column A and choosing Insert in Excel, then
# typing the new header and letting Excel auto-fit the column to the text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing columns (isin, mutualFundInstrumentName, planType, ...)
# one slot to the right, leaving a blank column A.
$ws.Columns.Item(1).Insert()

# Populate the new header cell.
$ws.Range("A1").Value = "owner_name"

# Auto-size the new column to fit its header text, like Excel would after
# double-clicking the column border.
$ws.Columns.Item(1).AutoFit()

# Leave the active cell on the new header, matching the post-edit selection.
[void]$ws.Range("A1").Select()
